# Adds two new defect rows (8 and 9) to the VoiceMaster_Defect sheet,
# mirroring the formatting of the existing data rows (e.g. row 22), and
# updates the sheet dimension / selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell text (kept as here-strings so embedded quotes / newlines
#     round-trip exactly without manual escaping) -------------------------

$textB23 = @"
User is able to overwrite user details by giving userId in :
curl -X POST \
  http://localhost:8080/api/customer 
"@

$textC23 = @"
IF we are placing this request: 
{
        "fname": "Madhu1",
        "lname": "Sharma",
        "email": "Madhu1",
        "secretCode": "9999",
        "userId":1
    }
It overwrites existing userId 1. Ideally it should not be allowed to enter UserId.

"@

$textB24 = "In SWAGGER,  example of above mentioned Request states UserID"

$textC24 = @"
This is the example of Add user API in swagger, which is incorrect as UserId should not be allowed. It should be generated automatically.
{
  "email": "string",
  "fname": "string",
  "lname": "string",
  "secretCode": "string",
  "userId": 0
}
"@

# --- Row 23 ---------------------------------------------------------------

$ws.Range("A22:D22").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A23").Value = 8
$ws.Range("B23").Value = $textB23
$ws.Range("C23").Value = $textC23
$ws.Rows("23:23").RowHeight = 165

# --- Row 24 ---------------------------------------------------------------

$ws.Range("A22:D22").Copy()
$ws.Range("A24:D24").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A24").Value = 9
$ws.Range("B24").Value = $textB24
$ws.Range("C24").Value = $textC24
$ws.Rows("24:24").RowHeight = 150

$excel.CutCopyMode = 0

# --- View state: put B23 as the active cell / selection -------------------

$ws.Range("B23").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
